# Slide 52 ("Computing Relative Addresses"), "Content Placeholder 2":
# the first paragraph's single run is split into three runs so the
# newly-typed word "for" becomes its own run, matching the author's
# incremental edit ("we did global variables" -> "we did for global
# variables").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(52)
$shape = $s.Shapes.Item("Content Placeholder 2")
$para = $shape.TextFrame.TextRange.Paragraphs(1)

# "we did " starts right after "Similar to what " (16 chars) and is
# 7 characters long; replacing just that slice with "we did for "
# reproduces the text selection the author retyped, leaving the runs
# before/after the edit untouched.
$sel = $para.Characters(17, 7)
$sel.Text = "we did for "
